# The sheet ships with cell protection enabled (sheetProtection), so the
# locked target cells can't be written to until the sheet is unprotected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# --- Update the "as of" date in the confidential disclosure banner (A10) ---
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-13 for illustrative purposes only and are subject to change."

# Setting a value with an embedded line break can make Excel auto-size the
# row; put the height back to the sheet default so nothing else shifts.
$ws.Rows(10).AutoFit()

# --- Refresh the Weight (D) and Percent Change (E) figures for rows 2-7 ---
$ws.Range("D2").Value = 0.2630257801799587
$ws.Range("E2").Value = -0.001372578923287904

$ws.Range("D3").Value = 0.5302104171010633
$ws.Range("E3").Value = -0.006033578174186593

$ws.Range("D4").Value = 0.05230293393759328
$ws.Range("E4").Value = -0.009271041628794774

$ws.Range("D5").Value = 0.09660007233371154
$ws.Range("E5").Value = -0.01529636711281079

$ws.Range("D6").Value = 0.05786079644767331
$ws.Range("E6").Value = -0.01810306530430927

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = -0.006570080266539158

# Restore sheet protection to match the original workbook's state.
$ws.Protect()
